$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add the new "Tom" comparison rows (82-105): additional longitudinal
# lipid-change comparisons for the CMT1A / PMP22 sciatic-nerve dataset,
# plus the missing metadata rows that go with them.
# ------------------------------------------------------------------

# Row 82-83: C22 vs C22-WT / C3 vs C3-WT  (experiment NLA_915)
$ws.Range("A82").Value2 = "NLA_915"
$ws.Range("B82").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 – C3)"
$ws.Range("C82").Value2 = "C22 vs C22 – WT"
$ws.Range("D82").Value2 = "NLA_099"

$ws.Range("A83").Value2 = "NLA_915"
$ws.Range("B83").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 – C3)"
$ws.Range("C83").Value2 = "C3 vs C3 – WT"
$ws.Range("D83").Value2 = "NLA_100"

# Row 84-85: same comparison, "no TG" variant (experiment NLA_916)
$ws.Range("A84").Value2 = "NLA_916"
$ws.Range("B84").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 – C3) (no TG)"
$ws.Range("C84").Value2 = "C22 vs C22 – WT"
$ws.Range("D84").Value2 = "NLA_101"

$ws.Range("A85").Value2 = "NLA_916"
$ws.Range("B85").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 – C3) (no TG)"
$ws.Range("C85").Value2 = "C3 vs C3 – WT"
$ws.Range("D85").Value2 = "NLA_102"

# Row 86-90: C22 age series (experiment NLA_917)
$ws.Range("A86").Value2 = "NLA_917"
$ws.Range("B86").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age)"
$ws.Range("C86").Value2 = "3 weeks"
$ws.Range("D86").Value2 = "NLA_103"

$ws.Range("A87").Value2 = "NLA_917"
$ws.Range("B87").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age)"
$ws.Range("C87").Value2 = "5 weeks"
$ws.Range("D87").Value2 = "NLA_104"

$ws.Range("A88").Value2 = "NLA_917"
$ws.Range("B88").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age)"
$ws.Range("C88").Value2 = "7 weeks"
$ws.Range("D88").Value2 = "NLA_105"

$ws.Range("A89").Value2 = "NLA_917"
$ws.Range("B89").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age)"
$ws.Range("C89").Value2 = "9 weeks"
$ws.Range("D89").Value2 = "NLA_106"

$ws.Range("A90").Value2 = "NLA_917"
$ws.Range("B90").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age)"
$ws.Range("C90").Value2 = "12 weeks"
$ws.Range("D90").Value2 = "NLA_107"

# Row 91-95: C22 age series, "no TG" variant (experiment NLA_918)
$ws.Range("A91").Value2 = "NLA_918"
$ws.Range("B91").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age) (no TG)"
$ws.Range("C91").Value2 = "3 weeks"
$ws.Range("D91").Value2 = "NLA_108"

$ws.Range("A92").Value2 = "NLA_918"
$ws.Range("B92").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age) (no TG)"
$ws.Range("C92").Value2 = "5 weeks"
$ws.Range("D92").Value2 = "NLA_109"

$ws.Range("A93").Value2 = "NLA_918"
$ws.Range("B93").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age) (no TG)"
$ws.Range("C93").Value2 = "7 weeks"
$ws.Range("D93").Value2 = "NLA_110"

$ws.Range("A94").Value2 = "NLA_918"
$ws.Range("B94").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age) (no TG)"
$ws.Range("C94").Value2 = "9 weeks"
$ws.Range("D94").Value2 = "NLA_111"

$ws.Range("A95").Value2 = "NLA_918"
$ws.Range("B95").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C22 age) (no TG)"
$ws.Range("C95").Value2 = "12 weeks"
$ws.Range("D95").Value2 = "NLA_112"

# Row 96-100: C3 age series (experiment NLA_919)
$ws.Range("A96").Value2 = "NLA_919"
$ws.Range("B96").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age)"
$ws.Range("C96").Value2 = "3 weeks"
$ws.Range("D96").Value2 = "NLA_113"

$ws.Range("A97").Value2 = "NLA_919"
$ws.Range("B97").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age)"
$ws.Range("C97").Value2 = "5 weeks"
$ws.Range("D97").Value2 = "NLA_114"

$ws.Range("A98").Value2 = "NLA_919"
$ws.Range("B98").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age)"
$ws.Range("C98").Value2 = "7 weeks"
$ws.Range("D98").Value2 = "NLA_115"

$ws.Range("A99").Value2 = "NLA_919"
$ws.Range("B99").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age)"
$ws.Range("C99").Value2 = "9 weeks"
$ws.Range("D99").Value2 = "NLA_116"

$ws.Range("A100").Value2 = "NLA_919"
$ws.Range("B100").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age)"
$ws.Range("C100").Value2 = "12 weeks"
$ws.Range("D100").Value2 = "NLA_117"

# Row 101-105: C3 age series, "no TG" variant (experiment NLA_920)
$ws.Range("A101").Value2 = "NLA_920"
$ws.Range("B101").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age) (no TG)"
$ws.Range("C101").Value2 = "3 weeks"
$ws.Range("D101").Value2 = "NLA_118"

$ws.Range("A102").Value2 = "NLA_920"
$ws.Range("B102").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age) (no TG)"
$ws.Range("C102").Value2 = "5 weeks"
$ws.Range("D102").Value2 = "NLA_119"

$ws.Range("A103").Value2 = "NLA_920"
$ws.Range("B103").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age) (no TG)"
$ws.Range("C103").Value2 = "7 weeks"
$ws.Range("D103").Value2 = "NLA_120"

$ws.Range("A104").Value2 = "NLA_920"
$ws.Range("B104").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age) (no TG)"
$ws.Range("C104").Value2 = "9 weeks"
$ws.Range("D104").Value2 = "NLA_121"

$ws.Range("A105").Value2 = "NLA_920"
$ws.Range("B105").Value2 = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A (C3 age) (no TG)"
$ws.Range("C105").Value2 = "12 weeks"
$ws.Range("D105").Value2 = "NLA_122"

# ------------------------------------------------------------------
# Formatting: the "group header" rows of each block (the rows that
# carry the merged per-experiment label, mirroring the shading used
# for every other multi-row experiment block in this sheet) get the
# same light-yellow fill used elsewhere (RGB FFF5CE); the remaining
# individual comparison rows stay unshaded, matching the rest of the
# sheet.
# ------------------------------------------------------------------
$highlightColor = 13563391  # RGB(255,245,206) -> FFF5CE, matches existing group rows

$ws.Range("A82:D83").Interior.Color = $highlightColor
$ws.Range("A86:D90").Interior.Color = $highlightColor
$ws.Range("A96:D100").Interior.Color = $highlightColor

$ws.Range("A84:D85").Interior.ColorIndex = -4142
$ws.Range("A91:D95").Interior.ColorIndex = -4142
$ws.Range("A101:D105").Interior.ColorIndex = -4142

# Rows 78-81 keep their existing (unshaded) look; just make sure no
# stray fill carried over so they stay visually identical.
$ws.Range("A78:D81").Interior.ColorIndex = -4142

# ------------------------------------------------------------------
# Column B grew wider to fit the new, longer comparison titles.
# (106.29 -> 110.46 "characters"; closest reachable value through the
# COM character-width rounding is used.)
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 109.6

# ------------------------------------------------------------------
# Selection / scroll position, mirroring where the author ended up
# after appending the new rows.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 62
$ws.Range("B101").Select()
